$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Natmi LR-pairs table rebuilt per Dr Hou's advice: the sending/target cluster
# cross-product now spans all three clusters (ECs, FAPs, sCs) instead of only
# FAPs/sCs x ECs/sCs, and every expression/specificity metric is recomputed.
$rows = @(
  @('ECs','Efnb3','Ephb1','ECs',1,0.3333333333333333,0.071327,0.213981,0.03356605248408491,0.03356605248408491,3,1,0.8299853333333332,2.489956,0.7350786001848651,0.7350786001848651,0.05920036387066666,0.532803274836,0.02467368687373285,0.02467368687373285),
  @('ECs','Efnb3','Ephb1','FAPs',1,0.3333333333333333,0.071327,0.213981,0.03356605248408491,0.03356605248408491,2,0.6666666666666666,0.01102233333333333,0.033067,0.009761957268446888,0.009761957268446888,0.0007861899696666667,0.007075709727,0.0003276703700200824,0.0003276703700200824),
  @('ECs','Efnb3','Ephb1','sCs',1,0.3333333333333333,0.071327,0.213981,0.03356605248408491,0.03356605248408491,3,1,0.2881033333333333,0.8643099999999999,0.2551594425466879,0.2551594425466879,0.02054954645666666,0.18494591811,0.008564695240331975,0.008564695240331974),
  @('FAPs','Efnb3','Ephb1','ECs',3,1,0.308223,0.9246690000000001,0.1450478695977975,0.1450478695977975,3,1,0.8299853333333332,2.489956,0.7350786001848651,0.7350786001848651,0.255820569396,2.302385124564,0.1066215849437459,0.1066215849437459),
  @('FAPs','Efnb3','Ephb1','FAPs',3,1,0.308223,0.9246690000000001,0.1450478695977975,0.1450478695977975,2,0.6666666666666666,0.01102233333333333,0.033067,0.009761957268446888,0.009761957268446888,0.003397336647,0.030576029823,0.001415951104892956,0.001415951104892956),
  @('FAPs','Efnb3','Ephb1','sCs',3,1,0.308223,0.9246690000000001,0.1450478695977975,0.1450478695977975,3,1,0.2881033333333333,0.8643099999999999,0.2551594425466879,0.2551594425466879,0.08880007371000001,0.79920066339,0.0370103335491587,0.03701033354915869),
  @('sCs','Efnb3','Ephb1','ECs',3,1,1.745424333333333,5.236273,0.8213860779181176,0.8213860779181176,3,1,0.8299853333333332,2.489956,0.7350786001848651,0.7350786001848651,1.448676597109777,13.038089373988,0.6037833283673865,0.6037833283673865),
  @('sCs','Efnb3','Ephb1','FAPs',3,1,1.745424333333333,5.236273,0.8213860779181176,0.8213860779181176,2,0.6666666666666666,0.01102233333333333,0.033067,0.009761957268446888,0.009761957268446888,0.01923864881011111,0.173147839291,0.008018335793533849,0.008018335793533849),
  @('sCs','Efnb3','Ephb1','sCs',3,1,1.745424333333333,5.236273,0.8213860779181176,0.8213860779181176,3,1,0.2881033333333333,0.8643099999999999,0.2551594425466879,0.2551594425466879,0.5028625685144444,4.525763116629999,0.2095844137571972,0.2095844137571972)
)

$r = 2
foreach ($row in $rows) {
  for ($c = 1; $c -le $row.Length; $c++) {
    $ws.Cells.Item($r, $c).Value = $row[$c - 1]
  }
  $r++
}
